$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Name" row (row 4) gets its value set: NiveauformacquisVs
$ws.Range("B4").Value = "NiveauformacquisVs"

# "Date" row (row 8) value updated to reflect new generation timestamp
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
